$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 612; this shifts all rows from 612..676 down to 613..677
# and extends the sheet dimension automatically to A1:R677.
$ws.Rows.Item(612).Insert()

# Populate the newly inserted row 612 with its data.
$ws.Cells.Item(612, 1).Value = 3
$ws.Cells.Item(612, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(612, 3).Value = "Coquimbo"
$ws.Cells.Item(612, 4).Value = 45194
$ws.Cells.Item(612, 5).Value = 5
$ws.Cells.Item(612, 6).Value = 100112031
$ws.Cells.Item(612, 7).Value = "Poroto verde"
$ws.Cells.Item(612, 8).Value = "Sin especificar"
$ws.Cells.Item(612, 9).Value = "Primera"
$ws.Cells.Item(612, 10).Value = 35
$ws.Cells.Item(612, 11).Value = 25000
$ws.Cells.Item(612, 12).Value = 25000
$ws.Cells.Item(612, 13).Value = 25000
$ws.Cells.Item(612, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(612, 15).Value = "Perú"
$ws.Cells.Item(612, 16).Value = 1000
$ws.Cells.Item(612, 17).Value = 25
$ws.Cells.Item(612, 18).Value = "Hortaliza"

# Apply the same date number format used by the other rows' Fecha (D)
# column to the newly inserted cell (matches cellXfs index 2 / numFmtId 165
# used throughout column D), rather than relying on Insert() to carry it
# over from row 611 above.
$ws.Cells.Item(612, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
